$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking price strings
# (e.g. "0.9989", "19.60") are preserved exactly as text, matching the
# original inlineStr cell type instead of being auto-converted to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.625.98'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.921.61'
$ws.Range('E3').Value = '  -0.77%  '
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '241.09'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = '0.4781'
$ws.Range('E7').Value = '  -2.09%  '
$ws.Range('D8').Value = '0.2867'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '0.06772'
$ws.Range('E9').Value = '  -1.59%  '
$ws.Range('D10').Value = '19.60'
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('D11').Value = '104.39'
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '0.07784'
$ws.Range('D13').Value = '1.935.16'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').Value = '5.258'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').Value = '0.6783'
$ws.Range('E15').Value = '  -3.06%  '
$ws.Range('D16').Value = '294.39'
$ws.Range('E16').Value = '  +7.86%  '
$ws.Range('D17').Value = '30.609.47'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '0.9989'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '0.000007560'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').Value = '12.83'
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').Value = '5.490'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '0.9989'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '6.373'
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '9.482'
$ws.Range('E24').Value = '  -3.43%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '168.32'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '19.71'
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = '2.116'
$ws.Range('E27').Value = '  -1.93%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '1.386'
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').Value = '0.1002'
$ws.Range('E29').Value = '  -3.47%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '4.594'
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.520'
$ws.Range('E31').Value = '  -2.32%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '4.301'
$ws.Range('E32').Value = '  -1.79%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.04774'
$ws.Range('E33').Value = '  -2.22%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '0.7331'
$ws.Range('E34').Value = '  -3.29%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.118'
$ws.Range('E35').Value = '  -2.64%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.711'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.01924'
$ws.Range('E37').Value = '  -4.20%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.627'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '6.392'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '75.20'
$ws.Range('E40').Value = '  -6.48%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '1.990'
$ws.Range('E41').Value = '  -4.38%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.8636'
$ws.Range('E42').Value = '  -4.37%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '106.38'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.4296'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').Value = '0.9992'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.485'
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '977.79'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '0.1210'
$ws.Range('E48').Value = '  -2.86%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '34.93'
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '8.880'
$ws.Range('E50').Value = '  -3.24%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.05816'
$ws.Range('E51').Value = '  +0.71%  '
